$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (Housing No.) for the new "Remarks" column
$ws.Columns("H").Insert()

# Match the new column's width to its neighbour (column G)
$ws.Columns("H").ColumnWidth = $ws.Columns("G").ColumnWidth

# Set header text for the new column
$ws.Range("H8").Value = "Remarks"

# Expand the merged title cell by one column to account for the inserted column
$ws.Range("B2:I2").UnMerge()
$ws.Range("B2:J2").Merge()

# Clear a couple of stray empty formatted cells left over from the old layout
$ws.Range("B1:D1").Clear()
$ws.Range("B7:D7").Clear()

# Update the active selection to reflect where the author left off editing
$ws.Range("B9").Select() | Out-Null
